# Applies the IFRS financial data corrections for rows 2-9 (company_list sheet)
# as described in the commit "error solve ifrs list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8106
$ws.Range("E2").Value = 307
$ws.Range("F2").Value = 307
$ws.Range("G2").Value = 538
$ws.Range("H2").Value = 382
$ws.Range("I2").Value = 404
$ws.Range("J2").Value = -22
$ws.Range("K2").Value = 8422
$ws.Range("L2").Value = 2013
$ws.Range("M2").Value = 6408
$ws.Range("N2").Value = 6289
$ws.Range("O2").Value = 119
$ws.Range("P2").Value = 521
$ws.Range("Q2").Value = 290
$ws.Range("R2").Value = -125
$ws.Range("S2").Value = -142
$ws.Range("T2").Value = 166
$ws.Range("U2").Value = 123
$ws.Range("V2").Value = 356
$ws.Range("W2").Value = 3.79
$ws.Range("X2").Value = 4.72
$ws.Range("Y2").Value = 6.38
$ws.Range("Z2").Value = 4.54
$ws.Range("AA2").Value = 31.42
$ws.Range("AB2").Value = 1113.89
$ws.Range("AC2").Value = 388
$ws.Range("AD2").Value = 18.13
$ws.Range("AE2").Value = 6929
$ws.Range("AF2").Value = 1.02
$ws.Range("AG2").Value = 230
$ws.Range("AH2").Value = 3.27
$ws.Range("AI2").Value = 52.03
$ws.Range("AJ2").Value = 84702850

# Row 3
$ws.Range("D3").Value = 8132
$ws.Range("E3").Value = 430
$ws.Range("F3").Value = 430
$ws.Range("G3").Value = 601
$ws.Range("H3").Value = 461
$ws.Range("I3").Value = 448
$ws.Range("J3").Value = 13
$ws.Range("K3").Value = 8206
$ws.Range("L3").Value = 1767
$ws.Range("M3").Value = 6440
$ws.Range("N3").Value = 6353
$ws.Range("O3").Value = 87
$ws.Range("P3").Value = 521
$ws.Range("Q3").Value = 523
$ws.Range("R3").Value = -169
$ws.Range("S3").Value = -419
$ws.Range("T3").Value = 174
$ws.Range("U3").Value = 349
$ws.Range("V3").Value = 283
$ws.Range("W3").Value = 5.29
$ws.Range("X3").Value = 5.66
$ws.Range("Y3").Value = 7.09
$ws.Range("Z3").Value = 5.54
$ws.Range("AA3").Value = 27.43
$ws.Range("AB3").Value = 1161.64
$ws.Range("AC3").Value = 430
$ws.Range("AD3").Value = 20.69
$ws.Range("AE3").Value = 7041
$ws.Range("AF3").Value = 1.26
$ws.Range("AG3").Value = 240
$ws.Range("AH3").Value = 2.7
$ws.Range("AI3").Value = 48.98
$ws.Range("AJ3").Value = 84702850

# Row 4
$ws.Range("D4").Value = 8207
$ws.Range("E4").Value = 428
$ws.Range("F4").Value = 428
$ws.Range("G4").Value = 536
$ws.Range("H4").Value = 418
$ws.Range("I4").Value = 426
$ws.Range("J4").Value = -8
$ws.Range("K4").Value = 8481
$ws.Range("L4").Value = 1853
$ws.Range("M4").Value = 6628
$ws.Range("N4").Value = 6548
$ws.Range("O4").Value = 80
$ws.Range("P4").Value = 521
$ws.Range("Q4").Value = 858
$ws.Range("R4").Value = -645
$ws.Range("S4").Value = -255
$ws.Range("T4").Value = 346
$ws.Range("U4").Value = 511
$ws.Range("V4").Value = 282
$ws.Range("W4").Value = 5.22
$ws.Range("X4").Value = 5.09
$ws.Range("Y4").Value = 6.6
$ws.Range("Z4").Value = 5.01
$ws.Range("AA4").Value = 27.96
$ws.Range("AB4").Value = 1198.79
$ws.Range("AC4").Value = 409
$ws.Range("AD4").Value = 19.89
$ws.Range("AE4").Value = 7266
$ws.Range("AF4").Value = 1.12
$ws.Range("AG4").Value = 240
$ws.Range("AH4").Value = 2.95
$ws.Range("AI4").Value = 51.21
$ws.Range("AJ4").Value = 84702850

# Row 5
$ws.Range("D5").Value = 8122
$ws.Range("E5").Value = 455
$ws.Range("F5").Value = 455
$ws.Range("G5").Value = 545
$ws.Range("H5").Value = 416
$ws.Range("I5").Value = 417
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 8571
$ws.Range("L5").Value = 1802
$ws.Range("M5").Value = 6769
$ws.Range("N5").Value = 6706
$ws.Range("O5").Value = 63
$ws.Range("P5").Value = 521
$ws.Range("Q5").Value = 929
$ws.Range("R5").Value = -1023
$ws.Range("S5").Value = -263
$ws.Range("T5").Value = 471
$ws.Range("U5").Value = 459
$ws.Range("V5").Value = 281
$ws.Range("W5").Value = 5.6
$ws.Range("X5").Value = 5.12
$ws.Range("Y5").Value = 6.3
$ws.Range("Z5").Value = 4.88
$ws.Range("AA5").Value = 26.61
$ws.Range("AB5").Value = 1238.71
$ws.Range("AC5").Value = 401
$ws.Range("AD5").Value = 20.46
$ws.Range("AE5").Value = 7501
$ws.Range("AF5").Value = 1.09
$ws.Range("AG5").Value = 240
$ws.Range("AH5").Value = 2.93
$ws.Range("AI5").Value = 51.8
$ws.Range("AJ5").Value = 84702850

# Row 6
$ws.Range("D6").Value = 7631
$ws.Range("E6").Value = 256
$ws.Range("F6").Value = 256
$ws.Range("G6").Value = 281
$ws.Range("H6").Value = 192
$ws.Range("I6").Value = 188
$ws.Range("K6").Value = 8303
$ws.Range("L6").Value = 2135
$ws.Range("M6").Value = 6168
$ws.Range("N6").Value = 6127
$ws.Range("P6").Value = 521
$ws.Range("Q6").Value = 642
$ws.Range("R6").Value = -250
$ws.Range("S6").Value = -414
$ws.Range("T6").Value = 222
$ws.Range("U6").Value = 419
$ws.Range("V6").Value = 299
$ws.Range("W6").Value = 3.36
$ws.Range("X6").Value = 2.52
$ws.Range("Y6").Value = 2.94
$ws.Range("Z6").Value = 2.28
$ws.Range("AA6").Value = 34.61
$ws.Range("AB6").Value = 1246.02
$ws.Range("AC6").Value = 181
$ws.Range("AD6").Value = 37.01
$ws.Range("AE6").Value = 7006
$ws.Range("AF6").Value = 0.96
$ws.Range("AG6").Value = 210
$ws.Range("AH6").Value = 3.13
$ws.Range("AI6").Value = 98.87
$ws.Range("AJ6").Value = 84702850

# Row 7
$ws.Range("D7").Value = 7571
$ws.Range("E7").Value = 290
$ws.Range("G7").Value = 300
$ws.Range("H7").Value = 230
$ws.Range("I7").Value = 222
$ws.Range("K7").Value = 8571
$ws.Range("L7").Value = 2528
$ws.Range("M7").Value = 6044
$ws.Range("N7").Value = 5998
$ws.Range("P7").Value = 520
$ws.Range("Q7").Value = 444
$ws.Range("R7").Value = -158
$ws.Range("S7").Value = -288
$ws.Range("T7").Value = 270
$ws.Range("W7").Value = 3.83
$ws.Range("X7").Value = 3.04
$ws.Range("Y7").Value = 3.65
$ws.Range("Z7").Value = 2.73
$ws.Range("AA7").Value = 41.82
$ws.Range("AC7").Value = 213
$ws.Range("AD7").Value = 26.7
$ws.Range("AE7").Value = 6975
$ws.Range("AF7").Value = 0.8100000000000001
$ws.Range("AG7").Value = 210
$ws.Range("AH7").Value = 3.7
$ws.Range("AI7").Value = 80.31
$ws.Range("U7").ClearContents()

# Row 8
$ws.Range("D8").Value = 7680
$ws.Range("E8").Value = 362
$ws.Range("G8").Value = 376
$ws.Range("H8").Value = 290
$ws.Range("I8").Value = 282
$ws.Range("K8").Value = 8649
$ws.Range("L8").Value = 2526
$ws.Range("M8").Value = 6128
$ws.Range("N8").Value = 6077
$ws.Range("P8").Value = 520
$ws.Range("Q8").Value = 739
$ws.Range("R8").Value = -480
$ws.Range("S8").Value = -208
$ws.Range("T8").Value = 443
$ws.Range("W8").Value = 4.71
$ws.Range("X8").Value = 3.78
$ws.Range("Y8").Value = 4.67
$ws.Range("Z8").Value = 3.37
$ws.Range("AA8").Value = 41.22
$ws.Range("AC8").Value = 271
$ws.Range("AD8").Value = 20.97
$ws.Range("AE8").Value = 7079
$ws.Range("AF8").Value = 0.8
$ws.Range("AG8").Value = 210
$ws.Range("AH8").Value = 3.7
$ws.Range("AI8").Value = 63.08
$ws.Range("U8").ClearContents()

# Row 9
$ws.Range("D9").Value = 7795
$ws.Range("E9").Value = 378
$ws.Range("G9").Value = 392
$ws.Range("H9").Value = 300
$ws.Range("I9").Value = 291
$ws.Range("K9").Value = 8742
$ws.Range("L9").Value = 2520
$ws.Range("M9").Value = 6222
$ws.Range("N9").Value = 6166
$ws.Range("P9").Value = 520
$ws.Range("Q9").Value = 758
$ws.Range("R9").Value = -476
$ws.Range("S9").Value = -207
$ws.Range("T9").Value = 443
$ws.Range("W9").Value = 4.84
$ws.Range("X9").Value = 3.84
$ws.Range("Y9").Value = 4.75
$ws.Range("Z9").Value = 3.44
$ws.Range("AA9").Value = 40.51
$ws.Range("AC9").Value = 279
$ws.Range("AD9").Value = 20.32
$ws.Range("AE9").Value = 7182
$ws.Range("AF9").Value = 0.79
$ws.Range("AG9").Value = 210
$ws.Range("AI9").Value = 61.13
$ws.Range("U9").ClearContents()
